$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update recalculated TPM-based ligand/receptor expression, specificity and
# edge-weight metrics for rows 2-10 (commit: "update scripts wuth new tpm").
$ws.Range("G2").Value = 6.240107999999999
$ws.Range("H2").Value = 18.720324
$ws.Range("I2").Value = 0.01732230523539376
$ws.Range("J2").Value = 0.01732230523539376
$ws.Range("M2").Value = 3.759736666666667
$ws.Range("N2").Value = 11.27921
$ws.Range("O2").Value = 0.0683751702595819
$ws.Range("P2").Value = 0.06837517025958188
$ws.Range("Q2").Value = 23.46116285156
$ws.Range("R2").Value = 211.15046566404
$ws.Range("S2").Value = 0.001184415569758495
$ws.Range("T2").Value = 0.001184415569758495
$ws.Range("G3").Value = 6.240107999999999
$ws.Range("H3").Value = 18.720324
$ws.Range("I3").Value = 0.01732230523539376
$ws.Range("J3").Value = 0.01732230523539376
$ws.Range("O3").Value = 0.6514180024294648
$ws.Range("P3").Value = 0.6514180024294647
$ws.Range("Q3").Value = 223.517159539268
$ws.Range("R3").Value = 2011.654435853412
$ws.Range("S3").Value = 0.01128406147391366
$ws.Range("T3").Value = 0.01128406147391366
$ws.Range("G4").Value = 6.240107999999999
$ws.Range("H4").Value = 18.720324
$ws.Range("I4").Value = 0.01732230523539376
$ws.Range("J4").Value = 0.01732230523539376
$ws.Range("O4").Value = 0.2802068273109533
$ws.Range("P4").Value = 0.2802068273109533
$ws.Range("Q4").Value = 96.14569123123999
$ws.Range("R4").Value = 865.3112210811599
$ws.Range("S4").Value = 0.004853828191721601
$ws.Range("T4").Value = 0.004853828191721601
$ws.Range("I5").Value = 0.9592798330716089
$ws.Range("J5").Value = 0.9592798330716091
$ws.Range("M5").Value = 3.759736666666667
$ws.Range("N5").Value = 11.27921
$ws.Range("O5").Value = 0.0683751702595819
$ws.Range("P5").Value = 0.06837517025958188
$ws.Range("Q5").Value = 1299.239337840863
$ws.Range("R5").Value = 11693.15404056777
$ws.Range("S5").Value = 0.06559092191285457
$ws.Range("T5").Value = 0.06559092191285457
$ws.Range("I6").Value = 0.9592798330716089
$ws.Range("J6").Value = 0.9592798330716091
$ws.Range("O6").Value = 0.6514180024294648
$ws.Range("P6").Value = 0.6514180024294647
$ws.Range("S6").Value = 0.6248921526303779
$ws.Range("T6").Value = 0.624892152630378
$ws.Range("I7").Value = 0.9592798330716089
$ws.Range("J7").Value = 0.9592798330716091
$ws.Range("O7").Value = 0.2802068273109533
$ws.Range("P7").Value = 0.2802068273109533
$ws.Range("S7").Value = 0.2687967585283764
$ws.Range("T7").Value = 0.2687967585283764
$ws.Range("G8").Value = 8.428738666666666
$ws.Range("I8").Value = 0.02339786169299727
$ws.Range("J8").Value = 0.02339786169299728
$ws.Range("M8").Value = 3.759736666666667
$ws.Range("N8").Value = 11.27921
$ws.Range("O8").Value = 0.0683751702595819
$ws.Range("P8").Value = 0.06837517025958188
$ws.Range("Q8").Value = 31.68983781881778
$ws.Range("R8").Value = 285.20854036936
$ws.Range("S8").Value = 0.001599832776968837
$ws.Range("T8").Value = 0.001599832776968838
$ws.Range("G9").Value = 8.428738666666666
$ws.Range("I9").Value = 0.02339786169299727
$ws.Range("J9").Value = 0.02339786169299728
$ws.Range("O9").Value = 0.6514180024294648
$ws.Range("P9").Value = 0.6514180024294647
$ws.Range("S9").Value = 0.01524178832517318
$ws.Range("T9").Value = 0.01524178832517318
$ws.Range("G10").Value = 8.428738666666666
$ws.Range("I10").Value = 0.02339786169299727
$ws.Range("J10").Value = 0.02339786169299728
$ws.Range("O10").Value = 0.2802068273109533
$ws.Range("P10").Value = 0.2802068273109533
$ws.Range("S10").Value = 0.006556240590855256
$ws.Range("T10").Value = 0.006556240590855256
